$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.054.76"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  +4.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.254.80"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  +5.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.49"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +2.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "619.74"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.12"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.388"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.248.51"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +5.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.797"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  -4.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.200"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.941.40"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +4.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "35.14"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.851.70"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  +5.03%  "
$ws.Range("B17").Value = "Toncoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.48"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.248.34"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  +4.93%  "
$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.63"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.59"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +16.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.02"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.86"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000202"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.10"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.58"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.13"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.89"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  -4.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.414.70"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +4.68%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.182"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("E31").Value = "  -5.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.21"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.11"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  +8.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.155"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.47"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -6.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.92"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.43"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "484.74"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.447"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.25"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.55"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -8.55%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.23"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "162.95"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.732"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +5.14%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.63"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +6.38%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.56"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.790"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  +8.71%  "
